# Actualización desde MV -datos-
# Append two new daily rows to the bottom of the data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newRows = @(
    @{ Fecha = "01-10-2021"; B = 13589; C = 19727; D = -6138 },
    @{ Fecha = "04-10-2021"; B = 13277; C = 19816; D = -6539 }
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

foreach ($row in $newRows) {
    $lastRow = $lastRow + 1
    $cellA = $ws.Cells.Item($lastRow, 1)
    # The date-like label (e.g. "01-10-2021") must land as literal text, not
    # get auto-converted into a date serial by Excel's input parser. Writing
    # it as a formula that evaluates to the text, then pasting-as-values,
    # stores a genuine shared-string cell with no extra number format/style.
    $cellA.Formula = '="' + $row.Fecha + '"'
    $cellA.Copy() | Out-Null
    $cellA.PasteSpecial(-4163) | Out-Null

    $ws.Cells.Item($lastRow, 2).Value = $row.B
    $ws.Cells.Item($lastRow, 3).Value = $row.C
    $ws.Cells.Item($lastRow, 4).Value = $row.D
}

$excel.CutCopyMode = 0
$ws.Range("A1").Select() | Out-Null
